$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.400.44"
$ws.Range("E2").Value = "  +0.19%  "
$ws.Range("D3").Value = "1.577.20"
$ws.Range("E3").Value = "  +0.15%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.73%  "
$ws.Range("E6").Value = "  +0.27%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.66"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.98%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "23.84"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.49%  "
$ws.Range("E10").Value = "  -0.02%  "
$ws.Range("E11").Value = "  -0.51%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0896"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.96%  "
$ws.Range("D13").Value = "1.802.89"
$ws.Range("E13").Value = "  +0.20%  "
$ws.Range("D14").Value = "1.572.79"
$ws.Range("E14").Value = "  -0.07%  "
$ws.Range("E15").Value = "  -0.07%  "
$ws.Range("D16").Value = "28.407.10"
$ws.Range("E16").Value = "  +0.39%  "
$ws.Range("E17").Value = "  -1.25%  "
$ws.Range("E18").Value = "  -0.80%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "230.16"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.33%  "
$ws.Range("E20").Value = "  +0.65%  "
$ws.Range("E21").Value = "  -0.77%  "
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.95"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.54%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.07"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.87%  "
$ws.Range("E25").Value = "  +2.26%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.52"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.35%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.98"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("E28").Value = "  -0.99%  "
$ws.Range("E29").Value = "  -0.64%  "
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("E31").Value = "  +3.86%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.09"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.19"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.18%  "
$ws.Range("E34").Value = "  -0.81%  "
$ws.Range("D35").Value = "1.394.46"
$ws.Range("E35").Value = "  +0.34%  "
$ws.Range("E36").Value = "  +7.79%  "
$ws.Range("E37").Value = "  -3.55%  "
$ws.Range("E38").Value = "  +0.20%  "
$ws.Range("E39").Value = "  +2.30%  "
$ws.Range("E40").Value = "  -0.32%  "
$ws.Range("E41").Value = "  -2.13%  "
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("E43").Value = "  +1.87%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.785"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.98%  "
$ws.Range("E45").Value = "  -3.28%  "
$ws.Range("E46").Value = "  -2.59%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.927"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.46%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "62.44"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.42%  "
$ws.Range("D49").Value = "1.714.92"
$ws.Range("E49").Value = "  +0.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "85.78"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.01%  "
$ws.Range("E51").Value = "  +3.74%  "
